$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.091.79"
$ws.Range("E2").Value = "  +1.77%  "

$ws.Range("D3").Value = "3.772.09"
$ws.Range("E3").Value = "  -0.52%  "

$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "623.53"
$ws.Range("E5").Value = "  +3.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.53"
$ws.Range("E6").Value = "  +1.41%  "

$ws.Range("D7").Value = "3.769.80"
$ws.Range("E7").Value = "  -0.54%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  +1.58%  "

$ws.Range("E10").Value = "  +1.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.460"
$ws.Range("E11").Value = "  +3.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.74"
$ws.Range("E12").Value = "  +1.41%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").Value = "  +0.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.72"
$ws.Range("E14").Value = "  +1.67%  "

$ws.Range("D15").Value = "4.411.14"
$ws.Range("E15").Value = "  -0.34%  "

$ws.Range("D16").Value = "3.772.34"
$ws.Range("E16").Value = "  -0.44%  "

$ws.Range("D17").Value = "69.132.62"
$ws.Range("E17").Value = "  +1.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.65"
$ws.Range("E18").Value = "  -2.73%  "

$ws.Range("E19").Value = "  -1.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.05"
$ws.Range("E20").Value = "  +0.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "467.33"
$ws.Range("E21").Value = "  +2.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.56"
$ws.Range("E22").Value = "  +1.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.706"
$ws.Range("E23").Value = "  +2.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000147"
$ws.Range("E24").Value = "  +2.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.16"
$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.05"
$ws.Range("E26").Value = "  +1.62%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.16"
$ws.Range("E27").Value = "  +3.89%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.03"
$ws.Range("E28").Value = "  +1.34%  "

$ws.Range("E29").Value = "  -0.14%  "

$ws.Range("D30").Value = "3.922.17"
$ws.Range("E30").Value = "  -0.37%  "

$ws.Range("E31").Value = "  +2.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.65"
$ws.Range("E32").Value = "  +2.63%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.16"
$ws.Range("E33").Value = "  -0.80%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.75"
$ws.Range("E34").Value = "  -0.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.169"
$ws.Range("E35").Value = "  +17.15%  "

$ws.Range("E36").Value = "  +0.18%  "

$ws.Range("D37").Value = "3.724.77"
$ws.Range("E37").Value = "  -0.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.95"
$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("E39").Value = "  +2.13%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.36"
$ws.Range("E40").Value = "  +5.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.82"
$ws.Range("E41").Value = "  +0.29%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.968"
$ws.Range("E42").Value = "  -1.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.13%  "

$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.32"
$ws.Range("E45").Value = "  -0.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.296"
$ws.Range("E46").Value = "  +0.53%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "152.33"
$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.68"
$ws.Range("E48").Value = "  -0.98%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.90"
$ws.Range("E49").Value = "  +3.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.40"
$ws.Range("E50").Value = "  +1.51%  "

$ws.Range("E51").Value = "  +0.92%  "
